$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set the text (shared-string) columns first, in the same order the
# original author entered them, so the shared-strings table is built in
# a matching sequence.
$ws.Range("A30").Value = "12월 9일"
$ws.Range("A31").Value = "12월 9일"
$ws.Range("A32").Value = "12월 10일"
$ws.Range("F30").Value = "index, show, update page jquery 작업"
$ws.Range("F31").Value = "footer, header jquery 작업, GUI Model 최종본 작성, 최종발표보고서 작성"
$ws.Range("F32").Value = "위키 작성, table page css/jquery, 최종산출물 보완 작업"

# Row 30: 12월 9일, 00:00 - 04:00, interruption 0, delta 240
$ws.Range("B30").Value = 0
$ws.Range("C30").Value = 0.16666666666666666
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 240

# Row 31: 12월 9일, 19:00 - 24:00, interruption 0, delta 300
$ws.Range("B31").Value = 0.79166666666666663
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 300

# Row 32: 12월 10일, 00:00 - 06:00, interruption 0, delta 360
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 0.25
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 360

# Update the view's scroll position and active selection to match the saved state
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("F32").Select()
